# Agrego modelo final con corte en 15.000
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: complete the "Modelo final" row that already had B11 set
$ws.Range("A11").Value = "SI"
$ws.Range("B11").Value = "Modelo final"
$ws.Range("C11").Value = "Entrenar el modelo final"
$ws.Range("D11").Value = "991_ZZ_lightgbm"
$ws.Range("E11").Value = "exp/HT9410/dataset_training.csv.gz"
$ws.Range("F11").Value = "exp/ZZ9410"

# Row 12: new row for the cut-at-15000 variant
$ws.Range("B12").Value = "Modelo final"
$ws.Range("C12").Value = "Entrenar el modelo final. Cortes hasta 15000"
$ws.Range("D12").Value = "991_ZZ_lightgbm_15000"
$ws.Range("E12").Value = "exp/HT9410/dataset_training.csv.gz"
$ws.Range("F12").Value = "exp/ZZ9411"

$ws.Range("B12").Select()

$ws.Columns.Item(4).ColumnWidth = 21.5
